$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @("18.71723686601473","9.378140749976851","8.879296903079922","33.22845121673124","3.66020160759573","24.44267138932052","10.11887976914525","12.07594715394128","25.17273522672529")
  3  = @("18.10766080884897","9.033147561697822","8.874897766933957","33.3606308283217","3.662363898944179","24.60128565113983","10.14852508564772","12.04690565584355","25.29278939868239")
  4  = @("17.72351906872316","8.813469573546193","8.87311349709333","33.45086546203495","3.663761981678636","24.70437569109774","10.16767817071315","12.03046191639189","25.37304514288324")
  5  = @("17.56471425146789","8.722072748622425","8.872618053399846","33.48991116020552","3.664349477498575","24.74781965883143","10.17572297383994","12.02411408146508","25.40739033538292")
  6  = @("17.53821475973935","8.706786137679442","8.872549808475897","33.49653180014337","3.664448105395698","24.75512012130366","10.17707330893115","12.02308147208175","25.41319226181835")
  7  = @("17.72138623962476","8.812244423240628","8.873105875943544","33.45138284777367","3.663769832845186","24.70495578488375","10.16778569401966","12.03037487220739","25.37350169881941")
  8  = @("18.50922473857615","9.260871786367749","8.877590485027257","33.27213842360291","3.660932582451596","24.49617883115706","10.12890456684368","12.06564814392332","25.21276929056405")
  9  = @("19.96743872948629","10.0745924863715","8.893608402008928","32.99299321429538","3.655924980076908","24.13197100240013","10.06017065780645","12.1456405951302","24.94969068006207")
  10 = @("20.97582833395969","10.6278784969929","8.909709575702211","32.83246233653096","3.652581390664777","23.89190922139273","10.0142064751467","12.21074071745573","24.78846310634822")
  11 = @("21.41912134041771","10.86918367157061","8.917958757072682","32.76919793354177","3.651132388646595","23.78866950052756","9.994271288181622","12.24166929777153","24.72213756847884")
  12 = @("21.58463985694081","10.95901612004477","8.921213932595089","32.74665200227309","3.650593986655613","23.75043315680571","9.986861717456659","12.25356453220755","24.69803593446469")
  13 = @("21.54909874722767","10.93973852608007","8.920507054614877","32.75144481049081","3.65070948372255","23.75862985599163","9.98845130889438","12.25099462380381","24.70318145321315")
  14 = @("21.43278636598173","10.87660546327873","8.918223938374309","32.76731474941815","3.651087887786297","23.78550656891364","9.993658907102017","12.24264428674914","24.72013435231699")
  15 = @("21.36123248338155","10.83773209987273","8.916842526488672","32.77721949583371","3.651321011505923","23.80208112291504","9.99686685353111","12.23755315557446","24.73065075380523")
  16 = @("20.94653631646676","10.61189510086559","8.909188945945282","32.83679387042798","3.652677531010536","23.8987762588227","10.01552883494113","12.20874542212698","24.79293930214587")
  17 = @("20.68808588178169","10.47065243671363","8.904729454600961","32.87584630053882","3.653528119185532","23.95962380717693","10.02722640958524","12.19140539984157","24.83295247783008")
  18 = @("20.53798649138988","10.38843823724658","8.902251607142714","32.89922665537593","3.654024136338892","23.9951831821679","10.03404627406207","12.18155593328321","24.85662722476486")
  19 = @("20.48692130457645","10.36043610722196","8.901427662793159","32.90730035827249","3.654193245509944","24.0073193975434","10.03637113680321","12.17824255103023","24.86475632031608")
  20 = @("20.71574895686672","10.48578927560134","8.905195168230588","32.8715939954254","3.653436871156607","23.95308836834227","10.02597169377073","12.19323847638209","24.82862463759587")
  21 = @("21.46701477500098","10.89519146424387","8.918890992017024","32.76261501625754","3.650976462114072","23.77758892864061","9.992125529580857","12.24509205639858","24.71512730245049")
  22 = @("21.94428156475577","11.15373767790964","8.928607060670613","32.69961796042028","3.649428476210904","23.66789301534094","9.970817648868845","12.28004670641869","24.64686537965068")
  23 = @("21.69085005946364","11.0165872084455","8.923351949217292","32.73248563177822","3.650249189365183","23.72598182145152","9.982115924637156","12.26129521742503","24.68275509025131")
  24 = @("20.70324718273678","10.47894905872996","8.904984351196337","32.87351357057662","3.653478102573687","23.95604124194503","10.02653865552778","12.19240936968737","24.83057916552845")
  25 = @("19.58334989288337","9.862009690349725","8.888508822487228","33.06071600656732","3.657220493876783","24.22566362935043","10.0779654566289","12.1456405951302","25.01525143208717")
}

$cols = @("B","C","D","F","G","I","J","L","O")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cellRef = "$($cols[$i])$row"
        $ws.Range($cellRef).Value = [double]$values[$i]
    }
}
